# feat: add 2022-Q3 data
#
# Target structure:
#   Sheets order: 总计, 2022-Q3 (new), 2022-Q2 (existing, unchanged content)
#   总计 sheet gains a new summary row for 2022-Q3 and keeps the old
#   2022-Q2 summary row (pushed down one row, with an incremented index).
#   A brand-new "2022-Q3" worksheet is populated with the fund holdings
#   table for that quarter.

$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# 1. Insert the new "2022-Q3" sheet right after "总计" (so the final
#    order is 总计 / 2022-Q3 / 2022-Q2).
# ------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $total)
$q3.Name = "2022-Q3"

# ------------------------------------------------------------------
# 2. Update the "总计" summary sheet.
#    Row 2 currently holds the 2022-Q2 summary (index 0). We push
#    that data down to row 3 (bumping its index to 1), then overwrite
#    row 2 with the new 2022-Q3 summary (index stays 0).
# ------------------------------------------------------------------

# Duplicate the formatting of row 2's index cell onto row 3 first, so
# the moved row keeps the same bold/bordered "index" style.
$total.Range("A2").Copy()
$total.Range("A3").PasteSpecial(-4122)

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0.67

# Now overwrite row 2 with the new quarter's summary values.
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 6
$total.Range("D2").Value = 1.31

# ------------------------------------------------------------------
# 3. Populate the new "2022-Q3" sheet with the fund holdings table.
# ------------------------------------------------------------------

# Reuse the existing bold/bordered "header" style (cell style index 2,
# already used on 总计!B1:D1 and 总计!A2) for the header row and the
# index column, instead of re-building the formatting from scratch.
$total.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$total.Range("A2").Copy()
$q3.Range("A2:A7").PasteSpecial(-4122)

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$col = 2
foreach ($h in $headers) {
    $q3.Cells.Item(1, $col).Value = $h
    $col = $col + 1
}

$rows = @(
    @("002345", "华夏高端制造灵活配置混合A", "20.58", "90.97", "3.51", "0.7224", 10),
    @("010490", "鹏华高质量增长混合A", "12.74", "93.98", "2.89", "0.3682", 10),
    @("009023", "鹏华稳健回报混合", "4.12", "94.39", "3.99", "0.1644", 8),
    @("015058", "华夏高端制造灵活配置混合C", "0.95", "90.97", "3.51", "0.0333", 10),
    @("010491", "鹏华高质量增长混合C", "0.44", "93.98", "2.89", "0.0127", 10),
    @("710002", "富安达策略精选混合", "0.59", "50.67", "2.10", "0.0124", 2)
)

# B:G hold text that looks numeric ("002345", "20.58", ...) — format the
# block as Text first so Excel stores them verbatim instead of coercing
# them into numbers (and dropping the leading zero on fund codes).
$q3.Range("B2:G7").NumberFormat = "@"

$r = 2
foreach ($row in $rows) {
    $q3.Cells.Item($r, 1).Value = $r - 2
    $q3.Cells.Item($r, 2).Value = $row[0]
    $q3.Cells.Item($r, 3).Value = $row[1]
    $q3.Cells.Item($r, 4).Value = $row[2]
    $q3.Cells.Item($r, 5).Value = $row[3]
    $q3.Cells.Item($r, 6).Value = $row[4]
    $q3.Cells.Item($r, 7).Value = $row[5]
    $q3.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

# The diff never touches <bookViews>, so leave the workbook's active
# sheet exactly where it started (总计, the first tab).
$total.Activate()
$total.Range("A1").Select() | Out-Null
